$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1493061553617529, 1.037716627549764, 2.445525955610853, 1.563817750126546, 1.572163572957211, 51)
    3  = @(-0.0001377995884945538, 0.8974840057777216, 1.83869274992964, 1.355984052240158, 1.369750733669702, 50)
    4  = @(0.1305475558742287, 0.9634510370349124, 2.686943648675822, 1.639189936729671, 1.650916085543009, 49)
    5  = @(0.02053268214733173, 0.8033519568783193, 1.685360403846863, 1.298214313527186, 1.311788311727695, 48)
    6  = @(0.1279059253806883, 1.069557746067552, 2.762993934606443, 1.662225596784757, 1.675214447118382, 47)
    7  = @(0.06123152309397346, 0.8594512018970224, 1.619385894554267, 1.272550939866168, 1.285122418519449, 46)
    8  = @(0.1270307523583763, 1.050670529041093, 2.656576264181475, 1.629900691508987, 1.643304426965771, 45)
    9  = @(0.06820288831713031, 0.9752875536895451, 2.031306037983135, 1.425238940663331, 1.440064547781096, 44)
    10 = @(0.1569399573777606, 1.109713135099426, 2.53526159535855, 1.592250481349762, 1.603249320694302, 43)
    11 = @(0.1971152062210202, 0.8222214692082306, 1.593335768711567, 1.262274046596684, 1.261901558681323, 42)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
